# Add a new experiment row (row 9) to Sheet1, documenting a run with
# bert-base-cased, no preprocessing/fine-tuning, dropout-rate tweak config,
# and its resulting metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$config = "seed = 1234`ntest_size = 0.2`nMAX_LEN = 512`nbatch_size = 16`nepochs = 10`nuse_gpu_predict = False`nlr = 1e-5"
$results = "accuracy    f1_macro    precision-neg    recall-neg`n----------  ----------  ---------------  ------------`n91.65%      91.65%      89.69%           94.09%"

$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "bert-base-cased"
$ws.Cells.Item(9, 3).Value = "NIL"
$ws.Cells.Item(9, 4).Value = "NIL"
$ws.Cells.Item(9, 5).Value = $config
$ws.Cells.Item(9, 6).Value = 0.94
$ws.Cells.Item(9, 7).Value = $results

$ws.Range("E9").WrapText = $true
$ws.Range("G9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 100.8

$ws.Range("B8").Select()
